# code formatting + show average and maximum droplet diameter

$wb = $excel.ActiveWorkbook
$plot = $wb.Worksheets.Item("Plot")
$droplets = $wb.Worksheets.Item("Droplet diameters")

# --- Bold the header rows (code formatting) ---
$plot.Range("A1:C1").Font.Bold = $true
$droplets.Range("A1").Font.Bold = $true

# --- Add average / maximum droplet diameter rows on the Plot sheet ---
$plot.Range("A7").Value = "Average droplet diameter (µm)"
$plot.Range("A7").Font.Bold = $true
$plot.Range("C7").Formula = "=AVERAGE('Droplet diameters'!A2:A12)"

$plot.Range("A8").Value = "Maximum droplet diameter (µm)"
$plot.Range("A8").Font.Bold = $true
$plot.Range("C8").Formula = "=MAX('Droplet diameters'!A2:A12)"

# --- Resize the histogram chart to make room for the new summary rows ---
$co = $plot.ChartObjects(1)
$co.Width = 561.9375
$co.Height = 291.75
